# Auto-generated edit script applying the row-pair swaps and odds updates
# described by the commit "Atualizacao de bases das ligas, do dia: 29-02-2024 as 07:50"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Range('B26').Value = 6810007
$ws.Range('F26').Value = 'Eupen'
$ws.Range('G26').Value = 'Club Brugge'
$ws.Range('I26').Value = 5
$ws.Range('K26').Value = 4.75
$ws.Range('L26').Value = 4
$ws.Range('M26').Value = 1.571
$ws.Range('N26').Value = 7
$ws.Range('O26').Value = 4.75
$ws.Range('P26').Value = 1.333
$ws.Range('Q26').Value = 1.5
$ws.Range('R26').Value = 1.875
$ws.Range('S26').Value = 1.975
$ws.Range('U26').Value = 1.925
$ws.Range('V26').Value = 1.925
$ws.Range('Y26').Value = 0.333
$ws.Range('AA26').Value = 0.9750000000000001
$ws.Range('AB26').Value = 0.925
$ws.Range('AC26').Value = -1

# Row 27
$ws.Range('B27').Value = 7030334
$ws.Range('F27').Value = 'Cercle Brugge'
$ws.Range('G27').Value = 'Genk'
$ws.Range('I27').Value = 1
$ws.Range('K27').Value = 2.75
$ws.Range('L27').Value = 3.5
$ws.Range('M27').Value = 2.25
$ws.Range('N27').Value = 2.4
$ws.Range('O27').Value = 3.5
$ws.Range('P27').Value = 2.55
$ws.Range('Q27').Value = 0
$ws.Range('R27').Value = 1.85
$ws.Range('S27').Value = 2
$ws.Range('U27').Value = 1.9
$ws.Range('V27').Value = 1.95
$ws.Range('Y27').Value = 1.55
$ws.Range('AA27').Value = 1
$ws.Range('AB27').Value = -1
$ws.Range('AC27').Value = 0.95

# Row 183
$ws.Range('B183').Value = 6810165
$ws.Range('F183').Value = 'Charleroi'
$ws.Range('G183').Value = 'Eupen'
$ws.Range('H183').Value = 1
$ws.Range('I183').Value = 0
$ws.Range('K183').Value = 1.6
$ws.Range('L183').Value = 4
$ws.Range('M183').Value = 5
$ws.Range('N183').Value = 1.8
$ws.Range('O183').Value = 3.75
$ws.Range('P183').Value = 4.2
$ws.Range('Q183').Value = -0.75
$ws.Range('R183').Value = 2.05
$ws.Range('S183').Value = 1.8
$ws.Range('T183').Value = 2.75
$ws.Range('U183').Value = 1.95
$ws.Range('V183').Value = 1.9
$ws.Range('W183').Value = 0.8
$ws.Range('Z183').Value = 0.5249999999999999
$ws.Range('AA183').Value = -0.5
$ws.Range('AB183').Value = -1
$ws.Range('AC183').Value = 0.8999999999999999

# Row 184
$ws.Range('B184').Value = 6810168
$ws.Range('F184').Value = 'OH Leuven'
$ws.Range('G184').Value = 'Genk'
$ws.Range('H184').Value = 2
$ws.Range('I184').Value = 1
$ws.Range('K184').Value = 4.5
$ws.Range('L184').Value = 4.2
$ws.Range('M184').Value = 1.666
$ws.Range('N184').Value = 4.333
$ws.Range('O184').Value = 4
$ws.Range('P184').Value = 1.7
$ws.Range('Q184').Value = 0.75
$ws.Range('R184').Value = 1.95
$ws.Range('S184').Value = 1.9
$ws.Range('T184').Value = 3
$ws.Range('U184').Value = 1.975
$ws.Range('V184').Value = 1.875
$ws.Range('W184').Value = 3.333
$ws.Range('Z184').Value = 0.95
$ws.Range('AA184').Value = -1
$ws.Range('AB184').Value = 0
$ws.Range('AC184').Value = 0

# Row 185
$ws.Range('B185').Value = 6810164
$ws.Range('F185').Value = 'Union Saint Gilloise'
$ws.Range('G185').Value = 'RWD Molenbeek'
$ws.Range('H185').Value = 3
$ws.Range('I185').Value = 2
$ws.Range('J185').Value = 'H'
$ws.Range('K185').Value = 1.2
$ws.Range('L185').Value = 7
$ws.Range('M185').Value = 12
$ws.Range('N185').Value = 1.142
$ws.Range('O185').Value = 8.5
$ws.Range('P185').Value = 15
$ws.Range('Q185').Value = -2.25
$ws.Range('R185').Value = 1.925
$ws.Range('S185').Value = 1.925
$ws.Range('T185').Value = 3.5
$ws.Range('U185').Value = 2.025
$ws.Range('V185').Value = 1.825
$ws.Range('W185').Value = 0.1419999999999999
$ws.Range('Y185').Value = -1
$ws.Range('AA185').Value = 0.925
$ws.Range('AB185').Value = 1.025
$ws.Range('AC185').Value = -1

# Row 186
$ws.Range('B186').Value = 6810162
$ws.Range('F186').Value = 'Standard Liege'
$ws.Range('G186').Value = 'Antwerp'
$ws.Range('H186').Value = 0
$ws.Range('I186').Value = 1
$ws.Range('J186').Value = 'A'
$ws.Range('K186').Value = 4
$ws.Range('L186').Value = 3.6
$ws.Range('M186').Value = 1.85
$ws.Range('N186').Value = 3.1
$ws.Range('O186').Value = 3.2
$ws.Range('P186').Value = 2.3
$ws.Range('Q186').Value = 0.25
$ws.Range('R186').Value = 1.8
$ws.Range('S186').Value = 2.05
$ws.Range('T186').Value = 2.25
$ws.Range('U186').Value = 1.875
$ws.Range('V186').Value = 1.975
$ws.Range('W186').Value = -1
$ws.Range('Y186').Value = 1.3
$ws.Range('AA186').Value = 1.05
$ws.Range('AB186').Value = -1
$ws.Range('AC186').Value = 0.9750000000000001

# Row 187
$ws.Range('B187').Value = 6810166
$ws.Range('F187').Value = 'KV Mechelen'
$ws.Range('G187').Value = 'Anderlecht'
$ws.Range('H187').Value = 2
$ws.Range('I187').Value = 2
$ws.Range('J187').Value = 'D'
$ws.Range('K187').Value = 3.5
$ws.Range('L187').Value = 3.5
$ws.Range('M187').Value = 2
$ws.Range('N187').Value = 3
$ws.Range('O187').Value = 3.5
$ws.Range('P187').Value = 2.2
$ws.Range('R187').Value = 1.925
$ws.Range('S187').Value = 1.925
$ws.Range('U187').Value = 1.875
$ws.Range('V187').Value = 1.975
$ws.Range('W187').Value = -1
$ws.Range('X187').Value = 2.5
$ws.Range('Z187').Value = 0.4625
$ws.Range('AA187').Value = -0.5
$ws.Range('AB187').Value = 0.875

# Row 188
$ws.Range('B188').Value = 6810163
$ws.Range('F188').Value = 'SintTruidense'
$ws.Range('G188').Value = 'Gent'
$ws.Range('H188').Value = 4
$ws.Range('I188').Value = 1
$ws.Range('J188').Value = 'H'
$ws.Range('K188').Value = 3.6
$ws.Range('L188').Value = 3.6
$ws.Range('M188').Value = 1.95
$ws.Range('N188').Value = 3.25
$ws.Range('O188').Value = 3.4
$ws.Range('P188').Value = 2.15
$ws.Range('R188').Value = 1.95
$ws.Range('S188').Value = 1.9
$ws.Range('U188').Value = 1.975
$ws.Range('V188').Value = 1.875
$ws.Range('W188').Value = 2.25
$ws.Range('X188').Value = 2.5
$ws.Range('Z188').Value = 0.95
$ws.Range('AA188').Value = -1
$ws.Range('AB188').Value = 0.9750000000000001

# Row 190
$ws.Range('B190').Value = 6810174
$ws.Range('F190').Value = 'Westerlo'
$ws.Range('G190').Value = 'OH Leuven'
$ws.Range('H190').Value = 0
$ws.Range('I190').Value = 3
$ws.Range('J190').Value = 'A'
$ws.Range('K190').Value = 1.909
$ws.Range('L190').Value = 3.75
$ws.Range('M190').Value = 3.5
$ws.Range('N190').Value = 1.909
$ws.Range('O190').Value = 3.5
$ws.Range('P190').Value = 3.8
$ws.Range('Q190').Value = -0.5
$ws.Range('R190').Value = 1.925
$ws.Range('S190').Value = 1.925
$ws.Range('U190').Value = 1.85
$ws.Range('V190').Value = 2
$ws.Range('W190').Value = -1
$ws.Range('Y190').Value = 2.8
$ws.Range('Z190').Value = -1
$ws.Range('AA190').Value = 0.925
$ws.Range('AB190').Value = 0.8500000000000001
$ws.Range('AC190').Value = -1

# Row 191
$ws.Range('B191').Value = 6810171
$ws.Range('F191').Value = 'KV Kortrijk'
$ws.Range('G191').Value = 'Charleroi'
$ws.Range('H191').Value = 1
$ws.Range('I191').Value = 0
$ws.Range('J191').Value = 'H'
$ws.Range('K191').Value = 3.2
$ws.Range('L191').Value = 3.5
$ws.Range('M191').Value = 2.1
$ws.Range('N191').Value = 3.4
$ws.Range('O191').Value = 3.4
$ws.Range('P191').Value = 2.05
$ws.Range('Q191').Value = 0.25
$ws.Range('R191').Value = 2
$ws.Range('S191').Value = 1.85
$ws.Range('U191').Value = 1.925
$ws.Range('V191').Value = 1.925
$ws.Range('W191').Value = 2.4
$ws.Range('Y191').Value = -1
$ws.Range('Z191').Value = 1
$ws.Range('AA191').Value = -1
$ws.Range('AB191').Value = -1
$ws.Range('AC191').Value = 0.925

# Row 221
$ws.Range('N221').Value = 2
$ws.Range('Q221').Value = -0.5
$ws.Range('R221').Value = 2.05
$ws.Range('S221').Value = 1.8

# Row 222
$ws.Range('N222').Value = 1.909
$ws.Range('P222').Value = 4
$ws.Range('R222').Value = 1.925
$ws.Range('S222').Value = 1.925

# Row 224
$ws.Range('N224').Value = 3.1
$ws.Range('O224').Value = 3.4
$ws.Range('P224').Value = 2.2
$ws.Range('U224').Value = 1.925
$ws.Range('V224').Value = 1.925

# Row 228
$ws.Range('N228').Value = 2.4
$ws.Range('P228').Value = 2.75
$ws.Range('U228').Value = 1.95
$ws.Range('V228').Value = 1.9
